# fix: bug in mapping to cognite types
#
# 1. The "source" property (ClassicEvent/ClassicAsset/ClassicFile rows) on the
#    Properties sheet pointed its Value Type at the non-existent
#    "cdf_cdm:SourceSystem(version=v1)" concept instead of the real
#    "cdf_cdm:CogniteSourceSystem(version=v1)" concept.
# 2. The "unitExternalId" property (ClassicTimeSeries row) on the Properties
#    sheet had a bare "Unit" placeholder instead of a proper reference to the
#    "cdf_cdm:CogniteUnit(version=v1)" concept.
# 3. Because (2) now references CogniteUnit, that view needs to actually be
#    defined on the Views sheet (alongside the already-present
#    CogniteSourceSystem view).

$wb = $excel.ActiveWorkbook

$propsSheet = $wb.Worksheets.Item("Properties")
$viewsSheet = $wb.Worksheets.Item("Views")

# --- Properties sheet --------------------------------------------------
# Fix "source" Value Type (rows 7, 14, 21 -> ClassicEvent, ClassicAsset, ClassicFile)
$propsSheet.Range("F7").Value2 = "cdf_cdm:CogniteSourceSystem(version=v1)"
$propsSheet.Range("F14").Value2 = "cdf_cdm:CogniteSourceSystem(version=v1)"
$propsSheet.Range("F21").Value2 = "cdf_cdm:CogniteSourceSystem(version=v1)"

# Fix "unitExternalId" Value Type (row 33 -> ClassicTimeSeries)
$propsSheet.Range("F33").Value2 = "cdf_cdm:CogniteUnit(version=v1)"

# --- Views sheet ---------------------------------------------------------
# Add the missing CogniteUnit view definition (row 9), right after the
# existing CogniteSourceSystem view (row 8).
$viewsSheet.Range("A9").Value2 = "cdf_cdm:CogniteUnit(version=v1)"
$viewsSheet.Range("C9").Value2 = "Represents a single unit of measurement"
$viewsSheet.Range("D9").Value2 = "CogniteDescribable"
$viewsSheet.Range("F9").Value2 = $true

$wb.Save()
